$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38/39: Coin name + Link swap (Frax <-> FraxShare), plus new price/volume data
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"

# Price (D) and Volume(1h) (E) updates
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.297.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.898.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.45%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4642"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3919"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07901"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9905"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.51%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.04"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.79%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.906.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.084"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.745"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06985"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.276.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.301"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.090"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.98%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.035"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "118.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.925"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.11%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.9063"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.302"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.329"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.231"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05817"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.178"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.65%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.000"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.806"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5702"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "

$ws.Range("E41").Value = "  -1.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.749"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.64%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5361"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07057"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.861"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.571"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "113.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.066"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.30"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "
